# Roll the FastReact Config tracker forward by one month.
#
# Every "From Date" (column C) and "End Date" (column D) for rows 2-14
# advances to the same point in the next calendar month: dates that sit
# on the first of the month move to the first of the next month, and
# dates that sit on the last day of the month move to the last day of
# the next month (rather than being clamped, e.g. Nov 30 -> Dec 31, not
# Nov 30 -> Dec 30).

function Add-OneMonthSmart {
    param($date)

    $firstOfThisMonth = $date.AddDays(1 - $date.Day)
    $firstOfNextMonth = $firstOfThisMonth.AddMonths(1)
    $lastOfThisMonth = $firstOfNextMonth.AddDays(-1)

    if ($date.Day -eq $lastOfThisMonth.Day) {
        # Original date was the last day of its month -> land on the
        # last day of the following month.
        $firstOfMonthAfterNext = $firstOfNextMonth.AddMonths(1)
        return $firstOfMonthAfterNext.AddDays(-1)
    } else {
        return $date.AddMonths(1)
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    foreach ($col in @('C', 'D')) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value()
        $cell.Value = Add-OneMonthSmart $current
    }
}
